# "actualizacion 13 de junio" — refresh the exported "archivo CVS" data:
# the FRONIUS inverter name/title changes, and the sampled dates move
# from May 2019 to scattered January 2019 dates (values stay the same).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header/title line — inverter model name updated
# (15.0-3 208 (1) (# 1)  ->  22.7-3 kW 1)
$ws.Range("A1").Value = "Fecha y hora,Energía por inversor|FRONIUS Symo 22.7-3 kW 1,Energía por inversor por kWp|FRONIUS Symo 22.7-3 kW 1,Instalación total,"

# Row 2: format-string row stays the same -> no change needed

# Rows 3-9: same daily figures, dates re-stamped to January 2019
$ws.Range("A3").Value = "01.01.2019,81.76,5.41,81.76,"
$ws.Range("A4").Value = "07.01.2019,70.85,4.69,70.85,"
$ws.Range("A5").Value = "09.01.2019,81.23,5.37,81.23,"
$ws.Range("A6").Value = "14.01.2019,76.61,5.07,76.61,"
$ws.Range("A7").Value = "16.01.2019,62.62,4.14,62.62,"
$ws.Range("A8").Value = "21.01.2019,63.90,4.23,63.90,"
$ws.Range("A9").Value = "30.01.2019,21.67,1.43,21.67,"

# Page setup now explicit (Letter size, portrait) for this sheet
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
